$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 20006358
$ws.Range("I74").Value = 50002996
$ws.Range("J74").Value = 8600
$ws.Range("K74").Value = 50002996
$ws.Range("L74").Value = 8600
$ws.Range("M74").Value = -50002060
$ws.Range("N74").Value = -10472

$ws.Range("H77").Value = 20006358
$ws.Range("I77").Value = 50002996
$ws.Range("J77").Value = 8600
$ws.Range("K77").Value = 250014980
$ws.Range("L77").Value = 43000
$ws.Range("M77").Value = -250010300
$ws.Range("N77").Value = -52360

$ws.Range("H100").Value = 16668266
$ws.Range("I100").Value = 16668266
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 16668266
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -16667725
$ws.Range("N100").ClearContents()

$ws.Range("H109").Value = 24061.975
$ws.Range("J109").Value = 24422.54
$ws.Range("L109").Value = 24422.54
$ws.Range("N109").Value = -27196.54

$ws.Range("H112").Value = 521026.34
$ws.Range("I112").Value = 685
$ws.Range("J112").Value = 550760.1
$ws.Range("K112").Value = 2055
$ws.Range("L112").Value = 1652280.3
$ws.Range("M112").Value = -947
$ws.Range("N112").Value = -1654496.3

$ws.Range("H125").Value = 906.0625
$ws.Range("J125").Value = 907.46155
$ws.Range("L125").Value = 8167.15395
$ws.Range("N125").Value = -13087.15395

$ws.Range("H132").Value = 563004.06
$ws.Range("I132").Value = 10374.182
$ws.Range("J132").Value = 1431422.4
$ws.Range("K132").Value = 31122.546
$ws.Range("L132").Value = 4294267.199999999
$ws.Range("M132").Value = -28592.546
$ws.Range("N132").Value = -4299327.199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2118.95
$ws.Range("I86").Value = 1935.2222
$ws.Range("J86").Value = 2269.2727
$ws.Range("K86").Value = 1935.2222
$ws.Range("L86").Value = 2269.2727
$ws.Range("M86").Value = -812.2221999999999
$ws.Range("N86").Value = -4515.2727

$ws.Range("H89").Value = 2118.95
$ws.Range("I89").Value = 1935.2222
$ws.Range("J89").Value = 2269.2727
$ws.Range("K89").Value = 9676.110999999999
$ws.Range("L89").Value = 11346.3635
$ws.Range("M89").Value = -4060.110999999999
$ws.Range("N89").Value = -22578.3635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2788.6956
$ws.Range("I31").Value = 1203.3636
$ws.Range("J31").Value = 4241.9165
$ws.Range("K31").Value = 1203.3636
$ws.Range("L31").Value = 4241.9165
$ws.Range("M31").Value = -908.3635999999999
$ws.Range("N31").Value = -4831.9165

$ws.Range("H34").Value = 2788.6956
$ws.Range("I34").Value = 1203.3636
$ws.Range("J34").Value = 4241.9165
$ws.Range("K34").Value = 1203.3636
$ws.Range("L34").Value = 4241.9165
$ws.Range("M34").Value = -1001.3636
$ws.Range("N34").Value = -4645.9165

$ws.Range("H58").Value = 1822.254
$ws.Range("I58").Value = 1713.3518
$ws.Range("J58").Value = 2475.6667
$ws.Range("K58").Value = 1713.3518
$ws.Range("L58").Value = 2475.6667
$ws.Range("M58").Value = -1510.3518
$ws.Range("N58").Value = -2881.6667

$ws.Range("H86").Value = 6077.2
$ws.Range("I86").Value = 5096.75
$ws.Range("K86").Value = 5096.75
$ws.Range("M86").Value = -3973.75

$ws.Range("H89").Value = 6077.2
$ws.Range("I89").Value = 5096.75
$ws.Range("K89").Value = 25483.75
$ws.Range("M89").Value = -19867.75

$ws.Range("H134").Value = 9681.538
$ws.Range("I134").Value = 11773.333
$ws.Range("J134").Value = 4975
$ws.Range("K134").Value = 35319.999
$ws.Range("L134").Value = 14925
$ws.Range("M134").Value = -32784.999
$ws.Range("N134").Value = -19995

$ws.Range("H136").Value = 1822.254
$ws.Range("I136").Value = 1713.3518
$ws.Range("J136").Value = 2475.6667
$ws.Range("K136").Value = 5140.055399999999
$ws.Range("L136").Value = 7427.000100000001
$ws.Range("M136").Value = -2590.055399999999
$ws.Range("N136").Value = -12527.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 15197.571
$ws.Range("I68").Value = 1131.3334
$ws.Range("J68").Value = 25747.25
$ws.Range("K68").Value = 3394.0002
$ws.Range("L68").Value = 77241.75
$ws.Range("M68").Value = -2583.0002
$ws.Range("N68").Value = -78863.75

$ws.Range("H71").Value = 15197.571
$ws.Range("I71").Value = 1131.3334
$ws.Range("J71").Value = 25747.25
$ws.Range("K71").Value = 10182.0006
$ws.Range("L71").Value = 231725.25
$ws.Range("M71").Value = -6126.000599999999
$ws.Range("N71").Value = -239837.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 27780590
$ws.Range("I80").Value = 35716972
$ws.Range("J80").Value = 3250
$ws.Range("K80").Value = 35716972
$ws.Range("L80").Value = 3250
$ws.Range("M80").Value = -35715974
$ws.Range("N80").Value = -5246

$ws.Range("H83").Value = 27780590
$ws.Range("I83").Value = 35716972
$ws.Range("J83").Value = 3250
$ws.Range("K83").Value = 178584860
$ws.Range("L83").Value = 16250
$ws.Range("M83").Value = -178579868
$ws.Range("N83").Value = -26234

$ws.Range("H97").Value = 741.86664
$ws.Range("I97").Value = 735.6667
$ws.Range("J97").Value = 766.6667
$ws.Range("K97").Value = 735.6667
$ws.Range("L97").Value = 766.6667
$ws.Range("M97").Value = -239.6667
$ws.Range("N97").Value = -1758.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2680.0833
$ws.Range("I7").Value = 1443.5294
$ws.Range("J7").Value = 5683.143
$ws.Range("K7").Value = 1443.5294
$ws.Range("L7").Value = 5683.143
$ws.Range("M7").Value = -1331.5294
$ws.Range("N7").Value = -5907.143

$ws.Range("H40").Value = 6352.0527
$ws.Range("I40").Value = 5399.231
$ws.Range("J40").Value = 8416.5
$ws.Range("K40").Value = 5399.231
$ws.Range("L40").Value = 8416.5
$ws.Range("M40").Value = -5263.231
$ws.Range("N40").Value = -8688.5

$ws.Range("H126").Value = 2680.0833
$ws.Range("I126").Value = 1443.5294
$ws.Range("J126").Value = 5683.143
$ws.Range("K126").Value = 4330.5882
$ws.Range("L126").Value = 17049.429
$ws.Range("M126").Value = -1860.5882
$ws.Range("N126").Value = -21989.429

$ws.Range("H132").Value = 4549.231
$ws.Range("I132").Value = 1130.4375
$ws.Range("J132").Value = 10019.3
$ws.Range("K132").Value = 3391.3125
$ws.Range("L132").Value = 30057.9
$ws.Range("M132").Value = -861.3125
$ws.Range("N132").Value = -35117.89999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7900
$ws.Range("I126").Value = 4350
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 13050
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -10580
$ws.Range("N126").Value = -49940

